# The document originally contains one paragraph split across several
# runs ("Aduuuuuuh" / " " / "si" / " asep") interleaved with spell-check
# <w:proofErr/> markers and a trailing "_GoBack" bookmark. The edit
# replaces all of that with a single clean run containing "bismillah".
#
# Simply overwriting the Range.Text leaves the proofErr/bookmark
# artifacts behind (they live outside the run/text stream), so instead
# we clear the whole story first - which drops that stray markup - and
# then insert the new text.

$d = $word.ActiveDocument

# Wipe all content (text + the leftover proofErr/bookmark markers).
$d.Content.Delete()

# Insert the replacement text into the now-empty document.
$d.Content.InsertAfter("bismillah")
